$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "66.014.77"
Set-TextValue $ws.Range("E2") "  +6.52%  "
Set-TextValue $ws.Range("D3") "3.019.42"
Set-TextValue $ws.Range("E3") "  +3.92%  "
Set-TextValue $ws.Range("E4") "  +0.01%  "
Set-TextValue $ws.Range("D5") "585.46"
Set-TextValue $ws.Range("E5") "  +2.84%  "
Set-TextValue $ws.Range("D6") "162.32"
Set-TextValue $ws.Range("E6") "  +12.85%  "
Set-TextValue $ws.Range("D7") "0.999"
Set-TextValue $ws.Range("E7") "  -0.06%  "
Set-TextValue $ws.Range("D8") "3.016.28"
Set-TextValue $ws.Range("E8") "  +3.89%  "
Set-TextValue $ws.Range("E9") "  +3.46%  "
Set-TextValue $ws.Range("D10") "6.76"
Set-TextValue $ws.Range("E10") "  -3.36%  "
Set-TextValue $ws.Range("E11") "  +5.78%  "
Set-TextValue $ws.Range("E12") "  +6.73%  "
Set-TextValue $ws.Range("D13") "0.0000255"
Set-TextValue $ws.Range("E13") "  +8.80%  "
Set-TextValue $ws.Range("D14") "34.79"
Set-TextValue $ws.Range("E14") "  +6.75%  "
Set-TextValue $ws.Range("E15") "  -0.53%  "
Set-TextValue $ws.Range("D16") "65.977.12"
Set-TextValue $ws.Range("E16") "  +6.56%  "
Set-TextValue $ws.Range("D17") "3.518.19"
Set-TextValue $ws.Range("E17") "  +3.87%  "
Set-TextValue $ws.Range("D18") "6.96"
Set-TextValue $ws.Range("E18") "  +7.14%  "
Set-TextValue $ws.Range("D19") "3.013.42"
Set-TextValue $ws.Range("E19") "  +3.78%  "
Set-TextValue $ws.Range("D20") "458.24"
Set-TextValue $ws.Range("E20") "  +6.65%  "
Set-TextValue $ws.Range("D21") "13.98"
Set-TextValue $ws.Range("E21") "  +7.46%  "
Set-TextValue $ws.Range("D22") "0.691"
Set-TextValue $ws.Range("E22") "  +6.13%  "
Set-TextValue $ws.Range("D23") "7.41"
Set-TextValue $ws.Range("E23") "  +7.99%  "
Set-TextValue $ws.Range("D24") "82.41"
Set-TextValue $ws.Range("E24") "  +4.63%  "
Set-TextValue $ws.Range("D25") "2.29"
Set-TextValue $ws.Range("E25") "  +12.87%  "
Set-TextValue $ws.Range("D26") "12.46"
Set-TextValue $ws.Range("E26") "  +3.98%  "
Set-TextValue $ws.Range("D27") "10.65"
Set-TextValue $ws.Range("E27") "  +3.95%  "
Set-TextValue $ws.Range("E28") "  -0.08%  "
Set-TextValue $ws.Range("D29") "8.11"
Set-TextValue $ws.Range("E29") "  +16.12%  "
Set-TextValue $ws.Range("D30") "2.36"
Set-TextValue $ws.Range("E30") "  +17.20%  "
Set-TextValue $ws.Range("D31") "0.0000107"
Set-TextValue $ws.Range("E31") "  -7.17%  "
Set-TextValue $ws.Range("E32") "  +4.31%  "
Set-TextValue $ws.Range("D33") "27.23"
Set-TextValue $ws.Range("E33") "  +6.42%  "
Set-TextValue $ws.Range("D34") "0.111"
Set-TextValue $ws.Range("E34") "  +4.26%  "
Set-TextValue $ws.Range("D35") "0.999"
Set-TextValue $ws.Range("E35") "  -0.08%  "
Set-TextValue $ws.Range("D36") "0.996"
Set-TextValue $ws.Range("E36") "  +4.42%  "
Set-TextValue $ws.Range("D37") "5.85"
Set-TextValue $ws.Range("E37") "  +8.51%  "
Set-TextValue $ws.Range("D38") "2.18"
Set-TextValue $ws.Range("E38") "  +15.04%  "
$ws.Range("B39").Value = "dogwifhat"
$ws.Range("C39").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue $ws.Range("D39") "3.00"
Set-TextValue $ws.Range("E39") "  +3.27%  "
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D40") "49.89"
Set-TextValue $ws.Range("E40") "  +2.21%  "
Set-TextValue $ws.Range("D41") "0.312"
Set-TextValue $ws.Range("E41") "  +17.25%  "
Set-TextValue $ws.Range("E42") "  +6.95%  "
Set-TextValue $ws.Range("D43") "43.77"
Set-TextValue $ws.Range("E43") "  +6.18%  "
Set-TextValue $ws.Range("D44") "8.46"
Set-TextValue $ws.Range("E44") "  +3.70%  "
Set-TextValue $ws.Range("D45") "395.65"
Set-TextValue $ws.Range("E45") "  +13.56%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D46") "0.0360"
Set-TextValue $ws.Range("E46") "  +6.97%  "
$ws.Range("B47").Value = "Maker"
$ws.Range("C47").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("D47") "2.801.09"
Set-TextValue $ws.Range("E47") "  +3.33%  "
Set-TextValue $ws.Range("D48") "133.86"
Set-TextValue $ws.Range("E48") "  +0.32%  "
Set-TextValue $ws.Range("D50") "23.89"
Set-TextValue $ws.Range("E50") "  +10.35%  "
Set-TextValue $ws.Range("E51") "  +4.41%  "
